$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-11 12:44:17"
$wsZh.Range("H2").Value = "2016-03-11 12:44:34"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-11 12:44:20"
$wsDe.Range("H2").Value = "2016-03-11 12:44:39"
